$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Places")

# Place indices for the new Cavetown entries
$ws.Range("A3").Value = 67
$ws.Range("A4").Value = 68
$ws.Range("A5").Value = 69
$ws.Range("A6").Value = 70

# Row 3 reuses the existing "Merchant" type string
$ws.Range("B3").Value = "Merchant"

# Remaining new cells, written in the order that reproduces the
# original author's shared-string insertion order
$ws.Range("B4").Value = "Blacksmith"
$ws.Range("B6").Value = "Raft Dealer"
$ws.Range("C3").Value = "Merchant Index 20: todo"
$ws.Range("C4").Value = "A bit more expensive (35) than Burnville blacksmith (25)"
$ws.Range("C5").Value = "A bit more expensive (20) than Burnville swim trainer (10)"
$ws.Range("B5").Value = "Swim Trainer"
$ws.Range("C6").Value = "350 Gold per raft, spawns east of Cavetown at the beach"

# Column B now needs to fit the new "Type" labels
$ws.Columns("B:B").ColumnWidth = 11.71

# Switch the active/selected sheet from GotoPoints to Places and move the
# selection to just past the newly added rows
[void]$ws.Activate()
[void]$ws.Range("C7").Select()
